# Updates cost_code (column G) values in the "Equip Billings" sheet.
#  - For every data row whose cost_code is exactly "9000 100F", append
#    " / CC NEEDED" to flag that a cost code is still required, EXCEPT
#    for rows 121-123 which instead get the combined
#    "9000 100M | 9000 100F | CC NEEDED" text.
#  - The trailing two "totals-style" helper rows (488 and 489) are
#    cleared out instead (row 488 fully blanked, row 489 only loses its
#    cost_code value) because they are stray placeholder rows, not real
#    billing lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 489
$specialRows = @(121, 122, 123)
$clearRows = @(488, 489)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($clearRows -contains $r) {
        continue
    }

    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "9000 100F") {
        if ($specialRows -contains $r) {
            $cell.Value = "9000 100M | 9000 100F | CC NEEDED"
        } else {
            $cell.Value = "9000 100F / CC NEEDED"
        }
    }
}

# Row 488 is a stray placeholder row: blank out its cost_code plus the
# numeric unit/rate/amount cells that held zeros.
$ws.Cells.Item(488, 7).Value = ""
$ws.Cells.Item(488, 9).Value = ""
$ws.Cells.Item(488, 11).Value = ""
$ws.Cells.Item(488, 12).Value = ""

# Row 489 keeps its numeric totals but loses the stray cost_code value.
$ws.Cells.Item(489, 7).Value = ""
